# Logged Week 16 and performed season sim from Week 17
$wb = $excel.ActiveWorkbook

# --- OFF sheet: update Home (row 2) target depth totals ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 399
$wsOff.Range("C2").Value = 285
$wsOff.Range("D2").Value = 111
$wsOff.Range("E2").Value = 59
$wsOff.Range("F2").Value = 8

# --- DEF sheet: update Home (row 2) target depth totals ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 386
$wsDef.Range("C2").Value = 270
$wsDef.Range("D2").Value = 106
$wsDef.Range("E2").Value = 55
$wsDef.Range("F2").Value = 9
$wsDef.Range("G2").Value = 4
